$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes the "Totals:" row and everything
# below it down by one, inheriting formatting from the row above).
$ws.Rows(4).Insert()

# Populate the newly inserted driver row (row 4) in the "Bad Drivers" table.
$ws.Range("A4").Value = "Fi - 16.0 (1657)"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 75
$ws.Range("D4").Value = 98.8

# Update the existing first driver row's Critical Minutes / Roaming % values.
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 98.7

# Update the "Totals:" row (now shifted down to row 5) to reflect the new
# driver that was added above it.
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 80

# Touch J20 formatting-only (no value) so the sheet's used range / dimension
# extends down to row 20 / column J, matching the extra trailing blank row
# that appears after the insert, without introducing any visible content.
$ws.Range("J20").Font.Bold = $false
